$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their original text formatting
# (values like "1.002" or "22.429.01" must not be auto-converted to numbers)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '22.429.01'
$ws.Range('E2').Value = '  +0.63%  '
$ws.Range('D3').Value = '1.572.14'
$ws.Range('E3').Value = '  +0.52%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '1.004'
$ws.Range('E5').Value = '  +0.27%  '
$ws.Range('D6').Value = '289.35'
$ws.Range('E6').Value = '  +0.17%  '
$ws.Range('D7').Value = '0.3745'
$ws.Range('E7').Value = '  -0.28%  '
$ws.Range('D8').Value = '49.02'
$ws.Range('E8').Value = '  -0.62%  '
$ws.Range('D9').Value = '0.3383'
$ws.Range('E9').Value = '  -0.81%  '
$ws.Range('B10').Value = 'Polygon'
$ws.Range('C10').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D10').Value = '1.127'
$ws.Range('E10').Value = '  -3.09%  '
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').Value = '0.07416'
$ws.Range('E11').Value = '  -2.94%  '
$ws.Range('D12').Value = '1.002'
$ws.Range('E12').Value = '  +0.10%  '
$ws.Range('D13').Value = '20.82'
$ws.Range('E13').Value = '  -2.46%  '
$ws.Range('D14').Value = '5.893'
$ws.Range('E14').Value = '  -1.83%  '
$ws.Range('D15').Value = '6.843'
$ws.Range('E15').Value = '  -1.22%  '
$ws.Range('D16').Value = '1.566.53'
$ws.Range('E16').Value = '  +0.46%  '
$ws.Range('D17').Value = '0.00001110'
$ws.Range('E17').Value = '  -1.51%  '
$ws.Range('D18').Value = '89.12'
$ws.Range('E18').Value = '  -0.87%  '
$ws.Range('D19').Value = '0.06684'
$ws.Range('E19').Value = '  -0.47%  '
$ws.Range('D20').Value = '1.005'
$ws.Range('E20').Value = '  +0.36%  '
$ws.Range('D21').Value = '6.130'
$ws.Range('E21').Value = '  -1.62%  '
$ws.Range('D22').Value = '16.10'
$ws.Range('E22').Value = '  -2.64%  '
$ws.Range('D23').Value = '11.77'
$ws.Range('E23').Value = '  -1.42%  '
$ws.Range('D24').Value = '22.416.32'
$ws.Range('E24').Value = '  +0.55%  '
$ws.Range('D25').Value = '2.358'
$ws.Range('E25').Value = '  -1.53%  '
$ws.Range('D26').Value = '2.540'
$ws.Range('E26').Value = '  -9.16%  '
$ws.Range('D27').Value = '19.95'
$ws.Range('E27').Value = '  -1.11%  '
$ws.Range('D28').Value = '147.22'
$ws.Range('E28').Value = '  +0.82%  '
$ws.Range('D29').Value = '4.981'
$ws.Range('E29').Value = '  +0.21%  '
$ws.Range('D30').Value = '124.82'
$ws.Range('E30').Value = '  -0.46%  '
$ws.Range('D31').Value = '1.740.80'
$ws.Range('E31').Value = '  +0.38%  '
$ws.Range('D32').Value = '1.982'
$ws.Range('E32').Value = '  -1.59%  '
$ws.Range('D33').Value = '0.9840'
$ws.Range('E33').Value = '  -3.65%  '
$ws.Range('D34').Value = '5.902'
$ws.Range('E34').Value = '  -4.74%  '
$ws.Range('D35').Value = '9.593'
$ws.Range('E35').Value = '  -4.29%  '
$ws.Range('B36').Value = 'Stellar'
$ws.Range('C36').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D36').Value = '0.08330'
$ws.Range('E36').Value = '  -2.36%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').Value = '1.384'
$ws.Range('E37').Value = '  +4.43%  '
$ws.Range('D38').Value = '0.02453'
$ws.Range('E38').Value = '  -3.22%  '
$ws.Range('D39').Value = '0.2243'
$ws.Range('E39').Value = '  -3.09%  '
$ws.Range('D40').Value = '0.06336'
$ws.Range('E40').Value = '  -0.82%  '
$ws.Range('D41').Value = '5.358'
$ws.Range('E41').Value = '  -2.62%  '
$ws.Range('D42').Value = '0.6189'
$ws.Range('E42').Value = '  -2.46%  '
$ws.Range('D43').Value = '11.03'
$ws.Range('E43').Value = '  -5.58%  '
$ws.Range('E44').Value = '  +0.50%  '
$ws.Range('D45').Value = '13.79'
$ws.Range('E45').Value = '  -2.62%  '
$ws.Range('D46').Value = '3.786'
$ws.Range('E46').Value = '  +1.01%  '
$ws.Range('D47').Value = '0.5765'
$ws.Range('E47').Value = '  -3.53%  '
$ws.Range('D48').Value = '2.040'
$ws.Range('E48').Value = '  -2.19%  '
$ws.Range('D49').Value = '125.38'
$ws.Range('E49').Value = '  +1.04%  '
$ws.Range('D50').Value = '1.222'
$ws.Range('E50').Value = '  -2.97%  '
$ws.Range('D51').Value = '0.07294'
$ws.Range('E51').Value = '  +0.55%  '
